# A&E lightning test cases
# Rebuild the "AdminUserDataQA" sheet (sheet4) into a 9-row x 7-column table
# of Draco automation test cases, and nudge the active-cell selection on the
# "LoginCredentials" sheet (sheet3).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet3 (LoginCredentials): only the remembered selection changes (C2 -> B3)
# Keep sheet4 as the active tab throughout, so only flip to sheet3 long
# enough to move its selection, then flip back.
# ---------------------------------------------------------------------------
$wsAdmin = $wb.Worksheets.Item("AdminUserDataQA")
$wsAdmin.Activate()

$wsLogin = $wb.Worksheets.Item("LoginCredentials")
$wsLogin.Range("B3").Select()

$wsAdmin.Activate()

# ---------------------------------------------------------------------------
# Sheet4 (AdminUserDataQA): new header row + 8 rows of Draco test-case data
# ---------------------------------------------------------------------------
$ws = $wsAdmin

# Headers
$ws.Cells.Item(1,1).Value = "Test Case ID"
$ws.Cells.Item(1,2).Value = "Environment"
$ws.Cells.Item(1,3).Value = "Salutation"
$ws.Cells.Item(1,4).Value = "First Name"
$ws.Cells.Item(1,5).Value = "Last Name"
$ws.Cells.Item(1,6).Value = "Email"
$ws.Cells.Item(1,7).Value = "Account Name"

# Data rows: TestCaseID, Environment, Salutation, FirstName, LastName, Email, AccountName
$rows = @(
    @("TC0001_DRACO","User1 QA","Mr.","Automation","DummyTest201","auto201@mailinator.com","Samsung Corp"),
    @("TC0002_DRACO","User1 QA","Mr.","Automation","DummyTest201","auto201@mailinator.com","Samsung Corp"),
    @("TC0004_DRACO","User1 QA","Mr.","Automation","DummyTest202","auto202@mailinator.com","Samsung Corp"),
    @("TC0003_DRACO","User1 QA","Mr.","Automation","DummyTest202","auto202@mailinator.com","Samsung Corp"),
    @("TC0005_DRACO","User1 QA","Mr.","Exportcheck","contact02","exportcontact2checktc5@mailinator.com","Automation_DRACO_Test"),
    @("TC0006_DRACO","User1 QA","Mr.","Exportcheck","contact01","exportcontactchecktc5@mailinator.com","Automation_DRACO_Test"),
    @("TC0007_DRACO","User1 QA","Mr.","employee8","pg8","employee8pg8@mailinator.com","season branch Ahemdabad"),
    @("TC0008_DRACO","User1 QA","Mr.","Testautomationuser01","contactdracodnt01","testautomationuser01contactdracodnt01@mailinator.com","season branch Ahemdabad")
)

$r = 2
foreach ($row in $rows) {
    $ws.Cells.Item($r,1).Value = $row[0]
    $ws.Cells.Item($r,2).Value = $row[1]
    $ws.Cells.Item($r,3).Value = $row[2]
    $ws.Cells.Item($r,4).Value = $row[3]
    $ws.Cells.Item($r,5).Value = $row[4]
    $ws.Cells.Item($r,6).Value = $row[5]
    $ws.Cells.Item($r,7).Value = $row[6]
    $r = $r + 1
}

# Hyperlink-ify the Email column (F2:F9), mailto: + the email text,
# matching the existing convention used on the other sheets.
for ($r = 2; $r -le 9; $r++) {
    $email = $ws.Cells.Item($r,6).Value
    $ws.Hyperlinks.Add($ws.Range("F" + $r), "mailto:" + $email)
    $ws.Range("F" + $r).Borders.LineStyle = 1
}

# ---------------------------------------------------------------------------
# Column widths for the new layout
# ---------------------------------------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 8.63
$ws.Columns.Item(3).ColumnWidth = 9.18
$ws.Columns.Item(4).ColumnWidth = 24.27
$ws.Columns.Item(5).ColumnWidth = 24.27
$ws.Columns.Item(6).ColumnWidth = 50.27
$ws.Columns.Item(7).ColumnWidth = 24.27

# ---------------------------------------------------------------------------
# View state: selection on the frozen-pane top-right pane moves to G9
# ---------------------------------------------------------------------------
$ws.Range("G9").Select()
